# se mejora las anulaciones y se modifica datos de sola la cuenta
$wb = $excel.ActiveWorkbook

# --- Update data on "DatosCuenta" sheet ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQAJuneTwo"
$wsCuenta.Range("B2").Value = "SmokeNameQAJuneTwo"
$wsCuenta.Range("C2").Value = 27100129
$wsCuenta.Range("D2").Value = 130

# Move selection on DatosCuenta to D10 and make it the active/selected sheet
$wsCuenta.Activate()
$wsCuenta.Range("D10").Select()

# --- "DatosAP" sheet keeps its own selection, just no longer the active tab ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("F14").Select()

# Re-activate DatosCuenta so it's the tab shown/selected when the file is saved
$wsCuenta.Activate()
